$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 130, shifting existing rows 130.. down by one.
$ws.Rows(130).Insert()

# Populate the newly inserted row 130 with the new data record.
$ws.Cells.Item(130, 1).Value = 5
$ws.Cells.Item(130, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(130, 3).Value = "Maule"
$ws.Cells.Item(130, 4).Value = 45090
$ws.Cells.Item(130, 5).Value = 7
$ws.Cells.Item(130, 6).Value = 100112031
$ws.Cells.Item(130, 7).Value = "Poroto verde"
$ws.Cells.Item(130, 8).Value = "Sin especificar"
$ws.Cells.Item(130, 9).Value = "Primera"
$ws.Cells.Item(130, 10).Value = 150
$ws.Cells.Item(130, 11).Value = 23000
$ws.Cells.Item(130, 12).Value = 23000
$ws.Cells.Item(130, 13).Value = 23000
$ws.Cells.Item(130, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(130, 15).Value = "Perú"
$ws.Cells.Item(130, 16).Value = 920
$ws.Cells.Item(130, 17).Value = 25
$ws.Cells.Item(130, 18).Value = "Hortaliza"
